# Applies the "Hoja De Datos" table update:
#  - Row 2 (Nery Javier de la Cruz Huinil) gets a spell-check proofErr
#    wrapper around the surname "Huinil".
#  - Row 3 (previously empty) gets "Hector Mauricio Cordero Oliva" /
#    "Desarrollador" (with the paragraph mark of the 2nd cell underlined).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wNs  = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$w14Ns = 'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- Row 2, Cell 1: split "Nery Javier de la Cruz Huinil" and mark
#     "Huinil" as a spell-check exception (w:proofErr spellStart/spellEnd).
$nameCell = $t.Cell(2, 1)
$nameParagraph = $nameCell.Range.Paragraphs.Item(1)
$nameXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="59DDD87F" w14:textId="3540C0CF" w:rsidR="00BF65F4" w:rsidRDefault="009D4EC1" w:rsidP="00BF65F4"><w:r><w:t xml:space="preserve">Nery Javier de la Cruz </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Huinil</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p>'
$nameParagraph.Range.InsertXML($nameXml)

# --- Row 3, Cell 1: fill in the developer's name.
$devNameCell = $t.Cell(3, 1)
$devNameParagraph = $devNameCell.Range.Paragraphs.Item(1)
$devNameXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="5046AA76" w14:textId="3192821C" w:rsidR="00BF65F4" w:rsidRDefault="00BF65F4" w:rsidP="00BF65F4"><w:r><w:t>Hector Mauricio Cordero Oliva</w:t></w:r></w:p>'
$devNameParagraph.Range.InsertXML($devNameXml)

# --- Row 3, Cell 2: fill in the role, underlining the paragraph mark.
$devRoleCell = $t.Cell(3, 2)
$devRoleParagraph = $devRoleCell.Range.Paragraphs.Item(1)
$devRoleXml = '<w:p ' + $wNs + ' ' + $w14Ns + ' w14:paraId="4206423D" w14:textId="5C9634C7" w:rsidR="00BF65F4" w:rsidRDefault="00BF65F4" w:rsidP="00BF65F4"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:t>Desarrollador</w:t></w:r></w:p>'
$devRoleParagraph.Range.InsertXML($devRoleXml)
